# Natmi following Dr Hou advice
#
# The sending-cluster set now also includes "sCs" (in addition to "FAPs"),
# so the Bmp7 -> Eng ligand-receptor table gains three more rows (sCs as
# sender, against the same three target clusters: ECs, FAPs, sCs) and the
# three original FAPs-sender rows get refreshed statistics to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (A1:T1 headers, unchanged):
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity,
# T Edge total expression derived specificity

function Set-NatmiRow($Row, $Sending, $Ligand, $Receptor, $Target, $Stats) {
    $ws.Cells.Item($Row, 1).Value = $Sending
    $ws.Cells.Item($Row, 2).Value = $Ligand
    $ws.Cells.Item($Row, 3).Value = $Receptor
    $ws.Cells.Item($Row, 4).Value = $Target

    for ($i = 0; $i -lt $Stats.Length; $i++) {
        $ws.Cells.Item($Row, 5 + $i).Value = $Stats[$i]
    }
}

# Rows 2-4: sending cluster "FAPs" (refreshed numbers)
Set-NatmiRow 2 "FAPs" "Bmp7" "Eng" "ECs" @(
    3, 1, 1.668521, 5.005563, 0.9677024783929865, 0.9677024783929865, 3, 1,
    170.93328, 512.79984, 0.7687311215213114, 0.7687311215213115,
    285.20576727888, 2566.85190550992, 0.7439030115139932, 0.7439030115139932
)

Set-NatmiRow 3 "FAPs" "Bmp7" "Eng" "FAPs" @(
    3, 1, 1.668521, 5.005563, 0.9677024783929865, 0.9677024783929865, 3, 1,
    40.31217066666667, 120.936512, 0.1812942463137967, 0.1812942463137967,
    67.26170331291733, 605.355329816256, 0.1754388914762496, 0.1754388914762496
)

Set-NatmiRow 4 "FAPs" "Bmp7" "Eng" "sCs" @(
    3, 1, 1.668521, 5.005563, 0.9677024783929865, 0.9677024783929865, 3, 1,
    11.112244, 33.336732, 0.04997463216489184, 0.04997463216489184,
    18.541012471124, 166.869112240116, 0.0483605754027437, 0.0483605754027437
)

# Rows 5-7: new sending cluster "sCs"
Set-NatmiRow 5 "sCs" "Bmp7" "Eng" "ECs" @(
    1, 0.3333333333333333, 0.05568766666666666, 0.167063, 0.03229752160701353, 0.03229752160701353, 3, 1,
    170.93328, 512.79984, 0.7687311215213114, 0.7687311215213115,
    9.51887551888, 85.66987966991999, 0.0248281100073183, 0.0248281100073183
)

Set-NatmiRow 6 "sCs" "Bmp7" "Eng" "FAPs" @(
    1, 0.3333333333333333, 0.05568766666666666, 0.167063, 0.03229752160701353, 0.03229752160701353, 3, 1,
    40.31217066666667, 120.936512, 0.1812942463137967, 0.1812942463137967,
    2.244890722695111, 20.204016504256, 0.005855354837547082, 0.005855354837547082
)

Set-NatmiRow 7 "sCs" "Bmp7" "Eng" "sCs" @(
    1, 0.3333333333333333, 0.05568766666666666, 0.167063, 0.03229752160701353, 0.03229752160701353, 3, 1,
    11.112244, 33.336732, 0.04997463216489184, 0.04997463216489184,
    0.6188149397906666, 5.569334458115999, 0.001614056762148148, 0.001614056762148148
)
